# Apply the "Mango pine apple sal nego" edit to the salary-comparison workbook.
#
# Target sheet for all the real content changes is the SECOND worksheet in the
# workbook ("Sheet3" internally / physical sheet2.xml) -- it is already the
# active/tabSelected sheet, matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- view state: zoom to 134% and move the selection to A9 -----------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 134
$ws.Range("A9").Select()

# --- A3/A4/A7: highlight part of the existing label text in red ------------
# (the wording itself is unchanged -- only a substring gets colored red)
$cellA3 = $ws.Range("A3")
$cellA3.Value = "Entire salary paid in monthly instalments (incl. holiday allowance and 13th month)"
$cellA3.Characters(67, 15).Font.Color = 255
$cellA3.Characters(82, 1).Font.Color = 0

$cellA4 = $ws.Range("A4")
$cellA4.Value = "Annual salary excl. Holiday allowance excl 13th month"
$cellA4.Characters(39, 15).Font.Color = 255

$cellA7 = $ws.Range("A7")
$cellA7.Value = "Annual salary included Holiday allowance and 13th month"
$cellA7.Characters(42, 14).Font.Color = 255

# --- new column H: "Net monthly" header + width ----------------------------
$ws.Columns.Item(8).ColumnWidth = 17.9
$ws.Range("H1").Value = "Net monthly"

# --- G2: newly filled-in net salary, formatted like its neighbour F2 -------
$ws.Range("G2").Value = 6000
$ws.Range("G2").NumberFormat = '[$€-2]\ #,##0;[Red]\-[$€-2]\ #,##0'

# --- H3: newly filled-in net monthly value, 2-decimal euro format ----------
$ws.Range("H3").Value = 4665.8
$ws.Range("H3").NumberFormat = '[$€-2]\ #,##0.00;[Red]\-[$€-2]\ #,##0.00'

# --- threaded-style review comments on F3, F4, F6 and G7 -------------------
$ws.Range("F3").AddComment("5500 + 440") | Out-Null
$ws.Range("F4").AddComment("5500 * 12") | Out-Null
$ws.Range("F6").AddComment("440 * 12") | Out-Null
$ws.Range("G7").AddComment("If you include 13th Month Pay.") | Out-Null

Write-Output "edit applied"
